$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.034145735084474
$ws.Range("D2").Value = 1.041121413211604
$ws.Range("E2").Value = 1.033367690110372
$ws.Range("F2").Value = 1.048459130576463
$ws.Range("I2").Value = 1.0375667191377
$ws.Range("J2").Value = 1.039266581152779
$ws.Range("K2").Value = 1.043901824984075
$ws.Range("L2").Value = 1.036170242187793
$ws.Range("M2").Value = 1.051218917063068
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.034984813698986
$ws.Range("D3").Value = 1.04177249187912
$ws.Range("E3").Value = 1.034078103663788
$ws.Range("F3").Value = 1.049292327645793
$ws.Range("I3").Value = 1.037749268715071
$ws.Range("J3").Value = 1.039749228151799
$ws.Range("K3").Value = 1.044363885799193
$ws.Range("L3").Value = 1.036689884545045
$ws.Range("M3").Value = 1.051864108976005
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.035528430544895
$ws.Range("D4").Value = 1.042194381532058
$ws.Range("E4").Value = 1.034538749150987
$ws.Range("F4").Value = 1.049832536580012
$ws.Range("I4").Value = 1.037866529989449
$ws.Range("J4").Value = 1.040061547934071
$ws.Range("K4").Value = 1.044662781548
$ws.Range("L4").Value = 1.037026414827872
$ws.Range("M4").Value = 1.052282016260929
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.035757127441987
$ws.Range("D5").Value = 1.042371885355157
$ws.Range("E5").Value = 1.034732632596459
$ws.Range("F5").Value = 1.050059895729283
$ws.Range("I5").Value = 1.037915619906845
$ws.Range("J5").Value = 1.04019284943908
$ws.Range("K5").Value = 1.044788414809556
$ws.Range("L5").Value = 1.037167959449952
$ws.Range("M5").Value = 1.052457804568791
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.035795536013313
$ws.Range("D6").Value = 1.042401697268221
$ws.Range("E6").Value = 1.034765199801361
$ws.Range("F6").Value = 1.05009808524513
$ws.Range("I6").Value = 1.03792385017694
$ws.Range("J6").Value = 1.04021489563481
$ws.Range("K6").Value = 1.044809507828392
$ws.Range("L6").Value = 1.037191729320088
$ws.Range("M6").Value = 1.052487326020141
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.035531485775948
$ws.Range("D7").Value = 1.042196752793438
$ws.Range("E7").Value = 1.034541338937182
$ws.Range("F7").Value = 1.049835573563365
$ws.Range("I7").Value = 1.037867186744496
$ws.Range("J7").Value = 1.040063302383921
$ws.Range("K7").Value = 1.044664460355905
$ws.Range("L7").Value = 1.037028305891389
$ws.Range("M7").Value = 1.05228436476171
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.034429164474249
$ws.Range("D8").Value = 1.041341323347627
$ws.Range("E8").Value = 1.033607577729449
$ws.Range("F8").Value = 1.048740489711904
$ws.Range("I8").Value = 1.037628590365844
$ws.Range("J8").Value = 1.039429690099738
$ws.Range("K8").Value = 1.044057998126678
$ws.Range("L8").Value = 1.036345797502649
$ws.Range("M8").Value = 1.051436873702178
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.032491991527006
$ws.Range("D9").Value = 1.039838608669183
$ws.Range("E9").Value = 1.031969608269257
$ws.Range("F9").Value = 1.046819130519203
$ws.Range("I9").Value = 1.037201594271539
$ws.Range("J9").Value = 1.038313355872042
$ws.Range("K9").Value = 1.04298871504944
$ws.Range("L9").Value = 1.03514539000193
$ws.Range("M9").Value = 1.049946814174888
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.031204177725443
$ws.Range("D10").Value = 1.0388400404941
$ws.Range("E10").Value = 1.030882738973657
$ws.Range("F10").Value = 1.045543935317518
$ws.Range("I10").Value = 1.036912563221057
$ws.Range("J10").Value = 1.037569325933638
$ws.Range("K10").Value = 1.042275522532921
$ws.Range("L10").Value = 1.034346721107687
$ws.Range("M10").Value = 1.048955777548082
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.030647423240081
$ws.Range("D11").Value = 1.038408441090866
$ws.Range("E11").Value = 1.030413347974076
$ws.Range("F11").Value = 1.04499314124028
$ws.Range("I11").Value = 1.036786382818681
$ws.Range("J11").Value = 1.037247214872848
$ws.Range("K11").Value = 1.041966639057712
$ws.Range("L11").Value = 1.034001286491819
$ws.Range("M11").Value = 1.048527222797917
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.030440753337595
$ws.Range("D12").Value = 1.038248245973003
$ws.Range("E12").Value = 1.030239182033229
$ws.Range("F12").Value = 1.044788760229751
$ws.Range("I12").Value = 1.036739360061041
$ws.Range("J12").Value = 1.037127578636641
$ws.Range("K12").Value = 1.041851897351541
$ws.Range("L12").Value = 1.033873037361508
$ws.Range("M12").Value = 1.048368125905095
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.030485078664856
$ws.Range("D13").Value = 1.038282602904351
$ws.Range("E13").Value = 1.030276532736336
$ws.Range("F13").Value = 1.044832591185682
$ws.Range("I13").Value = 1.036749453547648
$ws.Range("J13").Value = 1.037153240534779
$ws.Range("K13").Value = 1.041876510205935
$ws.Range("L13").Value = 1.033900544467078
$ws.Range("M13").Value = 1.048402248745483
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.030630337108294
$ws.Range("D14").Value = 1.038395196851407
$ws.Range("E14").Value = 1.030398947541945
$ws.Range("F14").Value = 1.044976242778484
$ws.Range("I14").Value = 1.036782499034933
$ws.Range("J14").Value = 1.037237325489486
$ws.Range("K14").Value = 1.041957154640755
$ws.Range("L14").Value = 1.033990684132857
$ws.Range("M14").Value = 1.048514070005958
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030719853370326
$ws.Range("D15").Value = 1.038464585670249
$ws.Range("E15").Value = 1.030474396150183
$ws.Range("F15").Value = 1.045064778946046
$ws.Range("I15").Value = 1.036802839099928
$ws.Range("J15").Value = 1.037289134390134
$ws.Range("K15").Value = 1.042006841235321
$ws.Range("L15").Value = 1.034046230237576
$ws.Range("M15").Value = 1.048582978410002
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.031241146316337
$ws.Range("D16").Value = 1.038868701068425
$ws.Range("E16").Value = 1.030913916997109
$ws.Range("F16").Value = 1.045580518825096
$ws.Range("I16").Value = 1.036920915799938
$ws.Range("J16").Value = 1.037590704735486
$ws.Range("K16").Value = 1.042296020827677
$ws.Range("L16").Value = 1.034369654894663
$ws.Range("M16").Value = 1.048984231482543
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.031568375699315
$ws.Range("D17").Value = 1.039122404044153
$ws.Range("E17").Value = 1.031189947611778
$ws.Range("F17").Value = 1.045904398077776
$ws.Range("I17").Value = 1.036994707460124
$ws.Range("J17").Value = 1.037779888521674
$ws.Range("K17").Value = 1.042477398716986
$ws.Range("L17").Value = 1.034572637208086
$ws.Range("M17").Value = 1.049236080905341
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.031759327604229
$ws.Range("D18").Value = 1.039270460527684
$ws.Range("E18").Value = 1.031351070205123
$ws.Range("F18").Value = 1.046093443787906
$ws.Range("I18").Value = 1.037037649647066
$ws.Range("J18").Value = 1.037890241809559
$ws.Range("K18").Value = 1.042583186873293
$ws.Range("L18").Value = 1.034691071304335
$ws.Range("M18").Value = 1.049383035358768
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.031824451557763
$ws.Range("D19").Value = 1.039320956768278
$ws.Range("E19").Value = 1.031406028905588
$ws.Range("F19").Value = 1.046157925889026
$ws.Range("I19").Value = 1.037052274968543
$ws.Range("J19").Value = 1.037927870336364
$ws.Range("K19").Value = 1.042619256757663
$ws.Range("L19").Value = 1.03473146066352
$ws.Range("M19").Value = 1.049433152298163
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.031533258326647
$ws.Range("D20").Value = 1.039095176259081
$ws.Range("E20").Value = 1.031160319857517
$ws.Range("F20").Value = 1.045869635187515
$ws.Range("I20").Value = 1.036986800573049
$ws.Range("J20").Value = 1.037759590304527
$ws.Range("K20").Value = 1.042457939248397
$ws.Range("L20").Value = 1.034550855191512
$ws.Range("M20").Value = 1.049209054150242
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.030587558418288
$ws.Range("D21").Value = 1.038362037409662
$ws.Range("E21").Value = 1.030362894261655
$ws.Range("F21").Value = 1.044933935194752
$ws.Range("I21").Value = 1.036772772200735
$ws.Range("J21").Value = 1.037212564275207
$ws.Range("K21").Value = 1.041933407092136
$ws.Range("L21").Value = 1.033964138561787
$ws.Range("M21").Value = 1.048481139004956
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.029993732594531
$ws.Range("D22").Value = 1.037901779793103
$ws.Range("E22").Value = 1.029862602705248
$ws.Range("F22").Value = 1.044346830408847
$ws.Range("I22").Value = 1.036637314942869
$ws.Range("J22").Value = 1.036868686698398
$ws.Range("K22").Value = 1.041603563640607
$ws.Range("L22").Value = 1.033595598234843
$ws.Range("M22").Value = 1.048023976499875
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.030308457041041
$ws.Range("D23").Value = 1.038145704486191
$ws.Range("E23").Value = 1.030127713572746
$ws.Range("F23").Value = 1.044657950683713
$ws.Range("I23").Value = 1.036709207437779
$ws.Range("J23").Value = 1.037050976637711
$ws.Range("K23").Value = 1.041778424115215
$ws.Range("L23").Value = 1.033790934620412
$ws.Range("M23").Value = 1.048266278406111
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.031549126094731
$ws.Range("D24").Value = 1.039107479089089
$ws.Range("E24").Value = 1.031173706998134
$ws.Range("F24").Value = 1.045885342632375
$ws.Range("I24").Value = 1.036990373661695
$ws.Range("J24").Value = 1.037768762178036
$ws.Range("K24").Value = 1.042466732164879
$ws.Range("L24").Value = 1.03456069743014
$ws.Range("M24").Value = 1.049221266207961
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.032992163812299
$ws.Range("D25").Value = 1.040226532950653
$ws.Range("E25").Value = 1.032392169956018
$ws.Range("F25").Value = 1.047314849980037
$ws.Range("I25").Value = 1.037312755847516
$ws.Range("J25").Value = 1.03860192670994
$ws.Range("K25").Value = 1.043265215573219
$ws.Range("L25").Value = 1.035455447736154
$ws.Range("M25").Value = 1.050331625711857
